$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 8")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# 1. "...complex scale-out products..." -> "...complex applications..."
$paraArch = $tr.Paragraphs(5)
$runArch = $paraArch.Runs(1)
$runArch.Text = "Experience in creating end-to-end solution architectures & designs for public cloud for complex applications across multiple domains including aspects such as cloud connectivity & networking, governance, security and enterprise IT integration."

# 2. "Building a collaborative..." -> "Built a collaborative..."
$paraBuilt = $tr.Paragraphs(7)
$runBuilt = $paraBuilt.Runs(1)
$runBuilt.Text = "Built a collaborative and transparent culture of continuous improvement, within the team and across domains, while mentoring junior team members."

# 3. Insert new bullet paragraph right after the "Built a collaborative..." paragraph.
$paraBuilt.InsertAfter("`rWorked in a coaching / collaborative style of working environment  and provided technical leadership.") | Out-Null
